# Fruta / hortaliza, semanal
# Update the weekly price records for "Feria Lagunitas de Puerto Montt - Tuna"
# Rows 2-7 of the sheet get their Fecha/Calidad/Volumen/Precio/Origen/Precio $/Kg
# values refreshed to the latest weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44250
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 806

# Row 3
$ws.Range("D3").Value = 44252
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 750

# Row 4
$ws.Range("D4").Value = 44257
$ws.Range("M4").Value = 100

# Row 5
$ws.Range("D5").Value = 44253
$ws.Range("M5").Value = 160

# Row 6
$ws.Range("D6").Value = 45072
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("R6").Value = "Provincia de Chacabuco"
$ws.Range("S6").Value = 889

# Row 7
$ws.Range("D7").Value = 45072
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 17000
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 944
